$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated probability matrix values from games pulled March 7
$ws.Range("B2").Value2 = 0.1993355481727575
$ws.Range("C2").Value2 = 0.553156146179402
$ws.Range("J2").Value2 = 0.009966777408637873
$ws.Range("P2").Value2 = 0.159468438538206
$ws.Range("S2").Value2 = 0.07807308970099668
$ws.Range("B3").Value2 = 0.005830903790087463
$ws.Range("C3").Value2 = 0.02040816326530612
$ws.Range("J3").Value2 = 0.04373177842565597
$ws.Range("P3").Value2 = 0.7142857142857143
$ws.Range("S3").Value2 = 0.2157434402332362
$ws.Range("J4").Value2 = 0.06756756756756757
$ws.Range("P4").Value2 = 0.6621621621621622
$ws.Range("S4").Value2 = 0.2702702702702703
$ws.Range("B6").Value2 = 0.06506024096385542
$ws.Range("D6").Value2 = 0.01204819277108434
$ws.Range("F6").Value2 = 0.06265060240963856
$ws.Range("J6").Value2 = 0.255421686746988
$ws.Range("O6").Value2 = 0.01445783132530121
$ws.Range("Q6").Value2 = 0.1253012048192771
$ws.Range("R6").Value2 = 0.07710843373493977
$ws.Range("S6").Value2 = 0.3879518072289156
$ws.Range("B7").Value2 = 0.07692307692307693
$ws.Range("D7").Value2 = 0.01758241758241758
$ws.Range("E7").Value2 = 0.004395604395604396
$ws.Range("F7").Value2 = 0.06373626373626373
$ws.Range("J7").Value2 = 0.1428571428571428
$ws.Range("O7").Value2 = 0.02637362637362637
$ws.Range("Q7").Value2 = 0.1494505494505494
$ws.Range("R7").Value2 = 0.08571428571428572
$ws.Range("S7").Value2 = 0.432967032967033
$ws.Range("B8").Value2 = 0.102820746132848
$ws.Range("D8").Value2 = 0.01910828025477707
$ws.Range("F8").Value2 = 0.05004549590536852
$ws.Range("J8").Value2 = 0.1091901728844404
$ws.Range("O8").Value2 = 0.02547770700636943
$ws.Range("Q8").Value2 = 0.1474067333939945
$ws.Range("R8").Value2 = 0.09372156505914468
$ws.Range("S8").Value2 = 0.4522292993630573
$ws.Range("B9").Value2 = 0.1068249258160237
$ws.Range("D9").Value2 = 0.01186943620178042
$ws.Range("F9").Value2 = 0.05044510385756677
$ws.Range("J9").Value2 = 0.09198813056379822
$ws.Range("O9").Value2 = 0.01483679525222552
$ws.Range("Q9").Value2 = 0.1750741839762611
$ws.Range("R9").Value2 = 0.09792284866468842
$ws.Range("S9").Value2 = 0.4510385756676558
$ws.Range("B10").Value2 = 0.1118038237738986
$ws.Range("D10").Value2 = 0.01620947630922693
$ws.Range("E10").Value2 = 0.0008312551953449709
$ws.Range("F10").Value2 = 0.07190357439733998
$ws.Range("J10").Value2 = 0.130091438071488
$ws.Range("O10").Value2 = 0.01537822111388196
$ws.Range("Q10").Value2 = 0.1978387364921031
$ws.Range("R10").Value2 = 0.08312551953449709
$ws.Range("S10").Value2 = 0.3728179551122194
$ws.Range("G11").Value2 = 0.1414790996784566
$ws.Range("J11").Value2 = 0.08520900321543408
$ws.Range("K11").Value2 = 0.1929260450160772
$ws.Range("L11").Value2 = 0.567524115755627
$ws.Range("S11").Value2 = 0.01286173633440514
$ws.Range("G12").Value2 = 0.7277628032345014
$ws.Range("J12").Value2 = 0.1752021563342318
$ws.Range("K12").Value2 = 0.01886792452830189
$ws.Range("L12").Value2 = 0.04582210242587601
$ws.Range("S12").Value2 = 0.03234501347708895
$ws.Range("G13").Value2 = 0.7659574468085106
$ws.Range("J13").Value2 = 0.2056737588652482
$ws.Range("S13").Value2 = 0.02836879432624113
$ws.Range("F15").Value2 = 0.009411764705882352
$ws.Range("H15").Value2 = 0.1858823529411765
$ws.Range("I15").Value2 = 0.05882352941176471
$ws.Range("J15").Value2 = 0.2776470588235294
$ws.Range("K15").Value2 = 0.07294117647058823
$ws.Range("M15").Value2 = 0.01411764705882353
$ws.Range("O15").Value2 = 0.08470588235294117
$ws.Range("S15").Value2 = 0.2964705882352941
$ws.Range("F16").Value2 = 0.01308900523560209
$ws.Range("H16").Value2 = 0.1910994764397906
$ws.Range("I16").Value2 = 0.06282722513089005
$ws.Range("J16").Value2 = 0.3952879581151832
$ws.Range("K16").Value2 = 0.1020942408376963
$ws.Range("M16").Value2 = 0.03926701570680628
$ws.Range("O16").Value2 = 0.05759162303664921
$ws.Range("S16").Value2 = 0.1387434554973822
$ws.Range("F17").Value2 = 0.01356350184956843
$ws.Range("H17").Value2 = 0.219482120838471
$ws.Range("I17").Value2 = 0.06658446362515413
$ws.Range("J17").Value2 = 0.3785450061652281
$ws.Range("K17").Value2 = 0.1159062885326757
$ws.Range("M17").Value2 = 0.02219482120838471
$ws.Range("N17").Value2 = 0.001233045622688039
$ws.Range("O17").Value2 = 0.05795314426633785
$ws.Range("S17").Value2 = 0.124537607891492
$ws.Range("F18").Value2 = 0.01470588235294118
$ws.Range("H18").Value2 = 0.2205882352941176
$ws.Range("I18").Value2 = 0.06127450980392157
$ws.Range("J18").Value2 = 0.4068627450980392
$ws.Range("K18").Value2 = 0.1053921568627451
$ws.Range("M18").Value2 = 0.03431372549019608
$ws.Range("O18").Value2 = 0.04411764705882353
$ws.Range("S18").Value2 = 0.1127450980392157
$ws.Range("F19").Value2 = 0.01530993278566094
$ws.Range("H19").Value2 = 0.2528005974607916
$ws.Range("I19").Value2 = 0.07729648991784914
$ws.Range("J19").Value2 = 0.3379387602688573
$ws.Range("K19").Value2 = 0.1064227035100822
$ws.Range("M19").Value2 = 0.03398058252427184
$ws.Range("N19").Value2 = 0.0007468259895444362
$ws.Range("O19").Value2 = 0.05974607916355489
$ws.Range("S19").Value2 = 0.1157580283793876
